$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "No movement" label to "No gesture" (cell content stays in A1; the
# shared-string table gets rebuilt on save so we just need the cell text right).
$ws.Range("A1").Value = "No gesture"

# Update the view: zoom to 190% and select A2 instead of H3.
$ws.Activate()
$excel.ActiveWindow.Zoom = 190
$ws.Range("A2").Select()
